$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -36.363636363636
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = -50

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -52.173913043478
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = -29.059829059829
$ws.Range("L16").Value = 22.058823529411
$ws.Range("M16").Value = -1.190476190476
$ws.Range("N16").Value = -83.976833976834

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 143
$ws.Range("J17").Value = 149
$ws.Range("K17").Value = -4.026845637583
$ws.Range("L17").Value = 26.548672566371
$ws.Range("M17").Value = 57.142857142857
$ws.Range("N17").Value = 0

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -18.75
$ws.Range("I18").Value = 109
$ws.Range("J18").Value = 106
$ws.Range("K18").Value = 2.830188679245
$ws.Range("L18").Value = 60.294117647058
$ws.Range("M18").Value = 91.228070175438
$ws.Range("N18").Value = -53.813559322033

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -36.363636363636
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = -39.130434782608
$ws.Range("I19").Value = 371
$ws.Range("J19").Value = 498
$ws.Range("K19").Value = -25.502008032128
$ws.Range("L19").Value = 15.576323987538
$ws.Range("M19").Value = 130.434782608696
$ws.Range("N19").Value = 29.268292682926

# --- Row 20 (G.L.A.) ---
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 38
$ws.Range("J20").Value = 47
$ws.Range("K20").Value = -19.148936170212
$ws.Range("L20").Value = -9.523809523809
$ws.Range("M20").Value = 8.571428571428
$ws.Range("N20").Value = -84.100418410041

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -34.146341463414
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = -32.394366197183
$ws.Range("I21").Value = 753
$ws.Range("J21").Value = 931
$ws.Range("K21").Value = -19.119226638023
$ws.Range("L21").Value = 21.256038647343
$ws.Range("M21").Value = 73.502304147465
$ws.Range("N21").Value = -47.853185595567

# --- Row 22 (Transit) ---
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 9.090909090909

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -38.461538461538
$ws.Range("I23").Value = 103
$ws.Range("J23").Value = 119
$ws.Range("K23").Value = -13.445378151260
$ws.Range("L23").Value = -11.206896551724
$ws.Range("M23").Value = 19.767441860465

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = -44
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 172
$ws.Range("H24").Value = -39.534883720930
$ws.Range("I24").Value = 792
$ws.Range("J24").Value = 1555
$ws.Range("K24").Value = -49.067524115755
$ws.Range("L24").Value = -6.603773584905
$ws.Range("M24").Value = 64.315352697095

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 23.333333333333
$ws.Range("I25").Value = 283
$ws.Range("J25").Value = 288
$ws.Range("K25").Value = -1.736111111111
$ws.Range("L25").Value = 58.988764044943
$ws.Range("M25").Value = 30.414746543778

# --- Row 26 (UCR Rape*) ---
$ws.Range("C26").Value = 1

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("L27").Value = 7.142857142857

# --- Row 28 (Shooting Vic.) ---
$ws.Range("G28").Value = "0"
$ws.Range("H28").Value = "***.*"

# --- Row 29 (Shooting Inc.) ---
$ws.Range("G29").Value = "0"
$ws.Range("H29").Value = "***.*"
